# Refresh the cryptocurrency price/volume snapshot in the worksheet.
# Source: diff of cryptos.xlsx scraped on 2024-04-19 (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several Price (column D) values are numeric-looking text (e.g. '556.07',
# '0.500') that must stay TEXT, exactly as stored in the source file, so trailing
# zeros/precision survive. A leading apostrophe forces Excel to keep them as text
# instead of auto-converting to a number.

$ws.Range('D2').Value = '63.764.17'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.054.38'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '''556.07'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  +3.24%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.050.13'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('D9').Value = '''0.500'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '''0.151'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Value = '''6.23'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').Value = '''34.94'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '3.570.62'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '63.770.25'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '3.056.15'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '''473.44'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('D21').Value = '''13.86'
$ws.Range('E21').Value = '  +1.70%  '
$ws.Range('D22').Value = '''0.671'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').Value = '''7.50'
$ws.Range('E23').Value = '  +3.87%  '
$ws.Range('D24').Value = '''13.43'
$ws.Range('E24').Value = '  +6.81%  '
$ws.Range('D25').Value = '''81.24'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').Value = '''2.78'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').Value = '''8.01'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '''2.05'
$ws.Range('E29').Value = '  +2.52%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = '''26.04'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('D36').Value = '''54.69'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').Value = '''458.63'
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('D38').Value = '''2.98'
$ws.Range('E38').Value = '  +13.53%  '
$ws.Range('D39').Value = '''0.0827'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').Value = '2.946.80'
$ws.Range('E41').Value = '  -6.58%  '
$ws.Range('D42').Value = '''8.25'
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('E43').Value = '  -5.29%  '
$ws.Range('D44').Value = '''27.87'
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').Value = '''0.258'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').Value = '''119.43'
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').Value = '0.0₃0513'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('E51').Value = '  -0.48%  '
